# Apply the two changes described by the commit:
#  1. Update the auto date/time footer placeholder text from 1/4/2024 to
#     1/6/2024 on the slide master and every slide layout.
#  2. Move the "Content Placeholder 2" text box on slide 2 to its new
#     position (size is unchanged).

$p = $ppt.ActivePresentation

# --- 1. Refresh the "Date Placeholder" footer field text everywhere ---
$oldDate = "1/4/2024"
$newDate = "1/6/2024"

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout belonging to the master
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# --- 2. Reposition the body text box on slide 2 ---
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item("Content Placeholder 2")
$contentShape.Left = 48
$contentShape.Top = 180
